# Update the "Förändrad" (Changed) date column (C) for all data rows
# from serial date 46061 (2026-02-08) to 46062 (2026-02-09).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$lastRow = $ws.Cells.Item($ws.Rows.Count, 3).End(-4162).Row

$ws.Range("C2:C$lastRow").Value = 46062
